$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Select()

# New shared random-data JSON template (replaces the two old static JSON
# bodies used for the POST and PUT request rows).
$newJson = "{`n  ""createdAt"": ""`$RandomPastDate"",`n  ""name"": ""`$RandomFullName"",`n  ""email"": ""`$RandomEmail"",`n  ""given_name"": ""`$RandomFirstName"",`n  ""last_ip"": ""`$RandomComputerIP"",`n  ""updated_at"": ""`$RandomPastDate"",`n  ""last_login"": ""`$RandomFutureDate"",`n  ""email_verified"": ""`$RandomBooleanValue""`n}"

$ws.Range("F2").Value = $newJson
$ws.Range("F5").Value = $newJson

# Row 5 shrinks to match row 2's height now that both rows hold the same
# (shorter) JSON template text.
$ws.Rows.Item(5).RowHeight = 145

# Scroll/selection state as left by the author after the edit.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F5").Select()
